$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 425; existing rows 425:492 shift down to 427:494
$ws.Rows("425:426").Insert()

# New row 425 (Primera)
$ws.Range("A425").Value = 1
$ws.Range("B425").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C425").Value = "Arica y Parinacota"
$ws.Range("D425").Value = 45142
$ws.Range("E425").Value = 15
$ws.Range("F425").Value = 100114014
$ws.Range("G425").Value = "Betarraga"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 1000
$ws.Range("K425").Value = 600
$ws.Range("L425").Value = 700
$ws.Range("M425").Value = 650
$ws.Range("N425").Value = "`$/paquete 4 unidades"
$ws.Range("O425").Value = "Región de Arica y Parinacota"
$ws.Range("P425").Value = 162
$ws.Range("Q425").Value = 4
$ws.Range("R425").Value = "Hortaliza"

# New row 426 (Segunda)
$ws.Range("A426").Value = 1
$ws.Range("B426").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C426").Value = "Arica y Parinacota"
$ws.Range("D426").Value = 45142
$ws.Range("E426").Value = 15
$ws.Range("F426").Value = 100114014
$ws.Range("G426").Value = "Betarraga"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Segunda"
$ws.Range("J426").Value = 1200
$ws.Range("K426").Value = 600
$ws.Range("L426").Value = 700
$ws.Range("M426").Value = 650
$ws.Range("N426").Value = "`$/paquete 5 unidades"
$ws.Range("O426").Value = "Región de Arica y Parinacota"
$ws.Range("P426").Value = 130
$ws.Range("Q426").Value = 5
$ws.Range("R426").Value = "Hortaliza"
